$wb = $excel.ActiveWorkbook

# --- Insert a new "Test_Cases" sheet in front of all existing sheets ---
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "Test_Cases"

# --- Header row ---
$ws.Range("A1").Value = "Sr.No"
$ws.Range("B1").Value = "TestCaseName"
$ws.Range("C1").Value = "Run_Mode"

# --- Data rows: test case name + run-mode flag ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "CreateAnNewAccount1"
$ws.Range("C2").Value = "Y"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "CreateAnNewAccount2"
$ws.Range("C3").Value = "Y"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "CreateAnNewAccount3"
$ws.Range("C4").Value = "N"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "TC04_Create_an_Account4"
$ws.Range("C5").Value = "N"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "TC05_Create_an_Account5"
$ws.Range("C6").Value = "N"

# --- Header formatting: bold, yellow fill, thin border, centered ---
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.Interior.Color = 65535
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.HorizontalAlignment = -4108

# --- Data formatting: thin border, centered ---
$data = $ws.Range("A2:C6")
$data.Borders.LineStyle = 1
$data.Borders.Weight = 2
$data.HorizontalAlignment = -4108

# --- Column widths to fit content ---
$ws.Columns.Item(1).ColumnWidth = 5.36328125
$ws.Columns.Item(2).ColumnWidth = 33.90625
$ws.Columns.Item(3).ColumnWidth = 10

$ws.Range("A1").Select()
